$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 87399297
$ws.Range("B4").Value = 77506
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = 'Garnlav'
$ws.Range("G4").Value = 'Alectoria sarmentosa'
$ws.Range("H4").Value = '(Ach.) Ach.'
$ws.Range("Q4").Value = 518535.8757363771
$ws.Range("R4").Value = 6965192.190086277

# Row 5
$ws.Range("A5").Value = 87397949
$ws.Range("B5").Value = 96239
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 504
$ws.Range("F5").Value = 'Guckusko'
$ws.Range("G5").Value = 'Cypripedium calceolus'
$ws.Range("H5").Value = 'L.'
$ws.Range("Q5").Value = 518333.776086097
$ws.Range("R5").Value = 6965141.215701581

# Row 6
$ws.Range("A6").Value = 101930649
$ws.Range("B6").Value = 96239
$ws.Range("E6").Value = 504
$ws.Range("F6").Value = 'Guckusko'
$ws.Range("G6").Value = 'Cypripedium calceolus'
$ws.Range("H6").Value = 'L.'
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = '5'
$ws.Range("I6").NumberFormat = "General"
$ws.Range("J6").Value = 'stjälkar/strån/skott'
$ws.Range("Q6").Value = 518336.0731258075
$ws.Range("R6").Value = 6965139.856789312
$ws.Range("S6").Value = 10
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = '2022-06-25'
$ws.Range("Y6").NumberFormat = "General"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = '2022-06-25'
$ws.Range("AA6").NumberFormat = "General"
$ws.Range("AC6").Value = '4 florala.'
$ws.Range("AW6").Value = 'Lars Grönvik'
$ws.Range("AX6").Value = 'Lars Grönvik'

# Row 7
$ws.Range("A7").Value = 87399184
$ws.Range("B7").Value = 96252
$ws.Range("E7").Value = 223591
$ws.Range("F7").Value = 'Skogsnycklar'
$ws.Range("G7").Value = 'Dactylorhiza maculata subsp. fuchsii'
$ws.Range("H7").Value = '(Druce) Hyl.'
$ws.Range("Q7").Value = 518353.1755671475
$ws.Range("R7").Value = 6965111.145664147

# Row 8
$ws.Range("A8").Value = 87399141
$ws.Range("B8").Value = 96254
$ws.Range("E8").Value = 223597
$ws.Range("F8").Value = 'Jungfru marie nycklar'
$ws.Range("G8").Value = 'Dactylorhiza maculata subsp. maculata'
$ws.Range("H8").Value = $null
$ws.Range("Q8").Value = 518219.8273314742
$ws.Range("R8").Value = 6965046.844427143

# Row 9
$ws.Range("A9").Value = 87399162
$ws.Range("B9").Value = 96356
$ws.Range("E9").Value = 219847
$ws.Range("F9").Value = 'Tvåblad'
$ws.Range("G9").Value = 'Neottia ovata'
$ws.Range("H9").Value = '(L.) Buff. & Fingerh.'
$ws.Range("Q9").Value = 518328.181063132
$ws.Range("R9").Value = 6965077.169283876

# Row 10
$ws.Range("A10").Value = 87399166
$ws.Range("B10").Value = 96356
$ws.Range("E10").Value = 219847
$ws.Range("F10").Value = 'Tvåblad'
$ws.Range("G10").Value = 'Neottia ovata'
$ws.Range("H10").Value = '(L.) Buff. & Fingerh.'
$ws.Range("Q10").Value = 518334.9957081943
$ws.Range("R10").Value = 6965086.809690522

# Row 11
$ws.Range("A11").Value = 87399150
$ws.Range("B11").Value = 97308
$ws.Range("E11").Value = 222467
$ws.Range("F11").Value = 'Gräsull'
$ws.Range("G11").Value = 'Eriophorum latifolium'
$ws.Range("H11").Value = 'Hoppe'

# Row 12
$ws.Range("A12").Value = 87399124
$ws.Range("B12").Value = 5135
$ws.Range("E12").Value = 105930
$ws.Range("F12").Value = 'Vågbandad barkbock'
$ws.Range("G12").Value = 'Semanotus undatus'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("Q12").Value = 518211.1325157742
$ws.Range("R12").Value = 6965045.881666865

# Row 13
$ws.Range("A13").Value = 87399156
$ws.Range("B13").Value = 96254
$ws.Range("E13").Value = 223597
$ws.Range("F13").Value = 'Jungfru marie nycklar'
$ws.Range("G13").Value = 'Dactylorhiza maculata subsp. maculata'
$ws.Range("H13").Value = $null
$ws.Range("I13").Value = $null
$ws.Range("J13").Value = $null
$ws.Range("Q13").Value = 518292.9341352677
$ws.Range("R13").Value = 6965075.143680119
$ws.Range("S13").Value = 25
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = '2020-08-09'
$ws.Range("Y13").NumberFormat = "General"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = '2020-08-09'
$ws.Range("AA13").NumberFormat = "General"
$ws.Range("AC13").Value = $null
$ws.Range("AW13").Value = 'Erland Lindblad'
$ws.Range("AX13").Value = 'Erland Lindblad'
